$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "done" column header
$ws.Range("H1").Value = "done"

# Set G column (header "blue") to "no" for all data rows, and fill the
# new H column ("done") with "no" as well.
for ($r = 2; $r -le 18; $r++) {
    $ws.Range("G$r").Value = "no"
    $ws.Range("H$r").Value = "no"
}

# Update the review text in F18 and wrap the text for that cell.
$ws.Range("F18").Value = "Take me to blockchain city!!!"
$ws.Range("F18").WrapText = $true

# Adjust the view: keep gridlines visible, scroll so column B is the
# top-left visible column, and select H2:H18 with H2 active.
$aw = $ws.Application.ActiveWindow
$aw.DisplayGridlines = $true
$aw.ScrollColumn = 2
$ws.Range("H2:H18").Select()
